$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1): update "想去人数" (column F) values
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value = 26
$ws1.Range("F6").Value = 260
$ws1.Range("F7").Value = 3641
$ws1.Range("F9").Value = 4257
$ws1.Range("F11").Value = 1063

# Sheet "全部类型" (sheet4): same events, different row numbers - update same column F values
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value = 26
$ws4.Range("F7").Value = 260
$ws4.Range("F8").Value = 3641
$ws4.Range("F10").Value = 4257
$ws4.Range("F12").Value = 1063
